# Add a "31/12/2023" data point for each region (Brasil, Nordeste, Sergipe),
# inserted right after that region's existing "31/12/2022" row, shifting the
# remaining rows down. The net effect is the data table growing from 16 rows
# (A1:C16) to 19 rows (A1:C19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data set (in final row order), header excluded.
$data = @(
  @("Brasil",   "31/12/2016", 4.667058984834714),
  @("Brasil",   "31/12/2017", 4.901081349408973),
  @("Brasil",   "31/12/2018", 4.844245820982904),
  @("Brasil",   "31/12/2019", 4.677684782843192),
  @("Brasil",   "31/12/2022", 4.695974113412442),
  @("Brasil",   "31/12/2023", 5.11059528529159),
  @("Nordeste", "31/12/2016", 3.520668787028191),
  @("Nordeste", "31/12/2017", 3.713599951282547),
  @("Nordeste", "31/12/2018", 3.897511354595247),
  @("Nordeste", "31/12/2019", 3.618295609666261),
  @("Nordeste", "31/12/2022", 3.483844743493894),
  @("Nordeste", "31/12/2023", 3.68496416112219),
  @("Sergipe",  "31/12/2016", 4.822422822579376),
  @("Sergipe",  "31/12/2017", 5.520020941397937),
  @("Sergipe",  "31/12/2018", 5.443832101077787),
  @("Sergipe",  "31/12/2019", 5.129675668138598),
  @("Sergipe",  "31/12/2022", 4.664592479082242),
  @("Sergipe",  "31/12/2023", 4.91480088324656)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value2 = $data[$i][0]
  $ws.Cells.Item($row, 2).Value2 = $data[$i][1]
  $ws.Cells.Item($row, 3).Value2 = $data[$i][2]
}
